$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '47.764.60'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.498.66'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.72%  '
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +1.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.09'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.14%  '
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.97'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.67%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.124'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.50%  '
$ws.Range("E14").Value = '  +0.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.889.20'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.494.77'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.849'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '47.657.03'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.26%  '
$ws.Range("E19").Value = '  +1.64%  '
$ws.Range("E20").Value = '  -0.49%  '
$ws.Range("E21").Value = '  +0.30%  '
$ws.Range("E22").Value = '  +12.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.75'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '247.66'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.07%  '
$ws.Range("E25").Value = '  -0.96%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.96'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.49%  '
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '35.03'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.35%  '
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.139'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.44%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.08'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.95'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.86'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.99%  '
$ws.Range("E34").Value = '  -2.17%  '
$ws.Range("E35").Value = '  -0.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("E37").Value = '  -1.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.66'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.96'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.111'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("E41").Value = '  +5.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '119.65'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.26%  '
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.003.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.81%  '
$ws.Range("E46").Value = '  +1.58%  '
$ws.Range("E47").Value = '  -3.49%  '
$ws.Range("E48").Value = '  +0.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("E50").Value = '  -2.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '56.87'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.11%  '
